# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets
# to reflect newly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 8416
    3  = 8001
    4  = 133
    5  = 194
    10 = 181
    11 = 239
    12 = 721
    13 = 145
    14 = 2068
    19 = 137
    20 = 46
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
